$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - new columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy header style from existing header (E1) to the new header cells
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Boolean outlier flag values for rows 2-18
$values = @(
    @($false, $false, $false),  # row 2
    @($false, $false, $false),  # row 3
    @($false, $false, $false),  # row 4
    @($false, $false, $false),  # row 5
    @($false, $false, $false),  # row 6
    @($false, $true,  $false),  # row 7
    @($false, $false, $false),  # row 8
    @($false, $false, $false),  # row 9
    @($false, $false, $false),  # row 10
    @($false, $false, $false),  # row 11
    @($false, $false, $false),  # row 12
    @($false, $true,  $false),  # row 13
    @($false, $false, $false),  # row 14
    @($false, $false, $false),  # row 15
    @($false, $false, $false),  # row 16
    @($false, $false, $false),  # row 17
    @($false, $true,  $false)   # row 18
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $rowVals = $values[$i]
    $ws.Cells.Item($row, 6).Value = $rowVals[0]
    $ws.Cells.Item($row, 7).Value = $rowVals[1]
    $ws.Cells.Item($row, 8).Value = $rowVals[2]
}
